$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hamburg")

# --- Reorder rows 10-12 ---
# Before:
#   row10: mask_date_if_bol_present = FALSE   (boolean, centered)
#   row11: g2_filename               = "G2 Schedule.xlsx"   (text)
#   row12: delay_filename            = "Vessel Delay Tracking.xlsx"   (text)
# After:
#   row10: g2_filename               = "G2 Schedule.xlsx"   (text)
#   row11: delay_filename            = "Vessel Delay Tracking.xlsx"   (text)
#   row12: mask_date_if_bol_present  = TRUE   (boolean, centered)

$ws.Cells.Item(10, 1).Value = "g2_filename"
$ws.Cells.Item(10, 2).Value = "G2 Schedule.xlsx"
$ws.Range("B10").HorizontalAlignment = 1          # xlHAlignGeneral - matches text rows (style index 7)

$ws.Cells.Item(11, 1).Value = "delay_filename"
$ws.Cells.Item(11, 2).Value = "Vessel Delay Tracking.xlsx"
$ws.Range("B11").HorizontalAlignment = 1          # xlHAlignGeneral - matches text rows (style index 7)

$ws.Cells.Item(12, 1).Value = "mask_date_if_bol_present"
$ws.Cells.Item(12, 2).Value = $true
$ws.Range("B12").HorizontalAlignment = -4108      # xlHAlignCenter - matches boolean rows (style index 9)

# --- Add new configuration row 17: g2_whitespace_rows = 9 ---
# Copy formatting from row 13 (a plain label/number row) so the new row
# picks up the same borders/alignment used throughout the table.
$ws.Range("A13:B13").Copy()
$ws.Range("A17:B17").PasteSpecial(-4122)          # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(17, 1).Value = "g2_whitespace_rows"
$ws.Cells.Item(17, 2).Value = 9

# --- Update the active cell selection shown when the sheet was last saved ---
$ws.Range("D6").Select()
